# Daily attendance processing - 2025-10-21 21:39:50
# For every data row, the "Recorded By" column (G) lists the users who
# recorded/modified the session, separated by ", ". Whenever the literal
# entry "System" appears in that list, promote it to the first position
# while keeping the relative order of the remaining entries untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ",\s*"

    # Use the .Equals() instance method (not -eq/-contains) so the
    # comparison stays case-sensitive: "System" must move, "system"
    # (lowercase) must stay where it is.
    $hasExactSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) { $hasExactSystem = $true }
    }

    if ($hasExactSystem) {
        $rest = @()
        foreach ($p in $parts) {
            if (-not $p.Equals("System")) { $rest += $p }
        }
        $newParts = @("System") + $rest
        $newVal = $newParts -join ", "

        if (-not $newVal.Equals($val)) {
            $cell.Value = $newVal
        }
    }
}
